$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "tabs": the "cells" tab now routes to a plain "cell" service
# instead of the old "PuffSmith\Cell\Import\CellImport" class name.
# ---------------------------------------------------------------------
$tabs = $wb.Worksheets.Item("tabs")
$tabs.Range("B2").Value = "cell"

# ---------------------------------------------------------------------
# Sheet "cells": convert the voltage (C) values from text "3.7" to a
# real number formatted with 2 decimals, mark the size (E) column as
# text-like (while keeping the already-entered numbers intact), and
# append a new row for the KeepPower IMR18350 cell.
# ---------------------------------------------------------------------
$cells = $wb.Worksheets.Item("cells")

# Voltage column: switch the number format to 2 decimals, then store
# the values as real numbers (order matters: while the cell is still
# text-formatted, assigning a number keeps it as text).
$cells.Range("C2:C5").NumberFormat = "0.00"
$cells.Range("C2").Value = 3.7
$cells.Range("C3").Value = 3.7
$cells.Range("C4").Value = 3.7
$cells.Range("C5").Value = 3.7

# Size column (E): the numbers are already stored in D2:D5 - re-assert
# them, then apply a text-like number format on top. Excel keeps
# existing numeric content untouched when a text format is applied
# after the fact, it only affects values typed in afterwards.
$cells.Range("E2").Value = 18650
$cells.Range("E3").Value = 18650
$cells.Range("E4").Value = 18650
$cells.Range("E5").Value = 21700
$cells.Range("E1:E5").NumberFormat = "@"

# New row 6: KeepPower IMR18350, 3.7V, 10A drain, 18350 size.
$cells.Range("A6").Value = "KeepPower"
$cells.Range("B6").Value = "IMR18350"
$cells.Range("C6").NumberFormat = "0.00"
$cells.Range("C6").Value = 3.7
$cells.Range("D6").Value = 10
$cells.Range("E6").Value = 18350
$cells.Range("E6").NumberFormat = "@"

# Keep the active selection near the newly added data, mirroring the
# author's last cursor position before saving.
$cells.Range("E4").Select()
